# Update "想去人数" (column F) values on the 展览 / 演出 / 全部类型 sheets
# to match freshly scraped counts (gh-pages data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 1409
$ws.Range("F6").Value = 307
$ws.Range("F7").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 728
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 152
$ws.Range("F15").Value = 109
$ws.Range("F16").Value = 0
$ws.Range("F18").Value = 189
$ws.Range("F20").Value = 394
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 55
$ws.Range("F27").Value = 0
$ws.Range("F29").Value = 15
$ws.Range("F30").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 558
$ws.Range("F33").Value = 56
$ws.Range("F34").Value = 2799
$ws.Range("F35").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("F39").Value = 1327
$ws.Range("F40").Value = 0
$ws.Range("F42").Value = 54
$ws.Range("F45").Value = 0
$ws.Range("F46").Value = 318

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 33

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 19781
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 7482
$ws.Range("F13").Value = 37
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("F18").Value = 189
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 394
$ws.Range("F21").Value = 71
$ws.Range("F25").Value = 64
$ws.Range("F26").Value = 317
$ws.Range("F28").Value = 0
$ws.Range("F30").Value = 175
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 3
$ws.Range("F36").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("F39").Value = 44
$ws.Range("F40").Value = 12565
$ws.Range("F43").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("F46").Value = 352
$ws.Range("F47").Value = 0
